$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (index 0): 99.98 -> 0M
$t.Cell(1, 1).Range.Text = "0M"

# Row 2 (index 1): 0.04 -> 0M
$t.Cell(2, 1).Range.Text = "0M"

# Row 3 (index 2): 231 -> 0M
$t.Cell(3, 1).Range.Text = "0M"

# Row 4 (index 3): 53 -> 212
$t.Cell(4, 1).Range.Text = "212"

# Row 6 (index 5): 0.00012 -> 0.00063
$t.Cell(6, 1).Range.Text = "0.00063"

# Row 7 (index 6): 0.00009 -> 0.00019
$t.Cell(7, 1).Range.Text = "0.00019"

# Row 8 (index 7): 0.00003 -> 0.00005
$t.Cell(8, 1).Range.Text = "0.00005"

# Row 9 (index 8): 0.00006 -> 0.00027
$t.Cell(9, 1).Range.Text = "0.00027"

# Row 10 (index 9): 0.00009 -> 0.00041
$t.Cell(10, 1).Range.Text = "0.00041"

# Row 11 (index 10): 0.00011 -> 0.00045
$t.Cell(11, 1).Range.Text = "0.00045"

# Row 12 (index 11): 0.00455 -> 0.04064
$t.Cell(12, 1).Range.Text = "0.04064"

# Row 44 (index 43): multi-value tab-separated -> 99.98
$t.Cell(44, 1).Range.Text = "99.98"

# Row 45 (index 44): multi-value tab-separated -> 0.04
$t.Cell(45, 1).Range.Text = "0.04"

# Row 46 (index 45): multi-value tab-separated -> 231
$t.Cell(46, 1).Range.Text = "231"
